$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 27 — this shifts the existing rows 27-29 down to 28-30,
# preserving their content/formatting, matching a weekly price-update where
# the newest week's record is inserted above the older ones.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new weekly entry.
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44783
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112043
$ws.Range("G27").Value = "Pepino dulce"
$ws.Range("H27").Value = "Cultivar IV Región"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17500
$ws.Range("N27").Value = "$/bandeja 18 kilos"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 972
$ws.Range("Q27").Value = 18
$ws.Range("R27").Value = "Hortaliza"
